$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 17 (test case "VALIDAZIONE_CDA2_LAB_CT7_KO", A17=53) -----------
# Timestamp / traceid / workflowInstanceId / execution date are no longer
# recorded for this case, and the applicability / rationale moves from
# "SI" + a "Referto prodotto..." note to "NO" with a new rationale.
$ws.Range("F17").ClearContents()
$ws.Range("G17").ClearContents()
$ws.Range("H17").ClearContents()
$ws.Range("I17").ClearContents()

$ws.Range("J17").Value = "NO"
$ws.Range("K17").Value = "L'applicativo è stato aggiornato in modo da salvare in automatico il CF in maiuscolo"

$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("O17").ClearContents()
$ws.Range("P17").ClearContents()

# --- Row 15 (test case "VALIDAZIONE_LAB_TIMEOUT", A15=44) ---------------
# Update the error-handling note to describe the new timeout behaviour.
$ws.Range("P15").Value = "In caso di timeout (errore 504) l'esecuzione prosegue ed il referto viene prodotto correttamente. Verrà visualizzato sulla richiesta l'esito negativo dell'invio a FSE e l'utente può effettuare il reinvio al momento oppure in un secondo momento"

# --- Selection / view state ---------------------------------------------
$ws.Range("P15").Select()
